$p = $ppt.ActivePresentation
try {
  $p.Slides.InsertFromFile("/tmp/work/before.pptx", 3)
  Write-Host ("Slides.Count=" + $p.Slides.Count)
} catch {
  Write-Host ("ERROR: " + $_)
}
